$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D1 header text: '外籍學者身分（教授、副教授、助理教授或博士後研究員）'
# -> '境外學者身分（教授、副教授、助理教授或博士後研究員）'
$ws.Range("D1").Value = "境外學者身分（教授、副教授、助理教授或博士後研究員）"

# Move the active selection to D1 (matches the saved selection state in the diff)
$ws.Range("D1").Select()
